$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 710
$ws.Cells.Item(3, 6).Value = 58
$ws.Cells.Item(4, 6).Value = 544
$ws.Cells.Item(5, 6).Value = 42
$bc = $ws.Cells.Item(7, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-15"
$bc.Style = "Normal"
$ws.Cells.Item(7, 3).Value = "蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）"
$ws.Cells.Item(7, 4).Value = "高新区望江西路888号 银泰百货（高新店）"
$ws.Cells.Item(7, 5).Value = "2024.09.15 10:00-10.02 22:00"
$ws.Cells.Item(7, 6).Value = 6
$ws.Cells.Item(7, 7).Value = 30
$ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91869"
$ws.Cells.Item(7, 9).Value = "//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png"
$bc = $ws.Cells.Item(8, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-16"
$bc.Style = "Normal"
$ws.Cells.Item(8, 3).Value = "肥西·星域动漫游戏嘉年华"
$ws.Cells.Item(8, 4).Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws.Cells.Item(8, 5).Value = "2024.09.16 10:00-09.16 17:00"
$ws.Cells.Item(8, 6).Value = 54
$ws.Cells.Item(8, 7).Value = 45
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws.Cells.Item(8, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"
$bc = $ws.Cells.Item(9, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-21"
$bc.Style = "Normal"
$ws.Cells.Item(9, 3).Value = "合肥·漫有引力动漫游戏嘉年华"
$ws.Cells.Item(9, 4).Value = "幸福路1号(筑梦集团·结婚产业园·B1幢) 费加罗宴会艺术中心(旗舰店)"
$ws.Cells.Item(9, 5).Value = "2024.09.21 10:00-09.21 17:00"
$ws.Cells.Item(9, 6).Value = 48
$ws.Cells.Item(9, 7).Value = 50
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90905"
$ws.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/gfeOndjM1723659151069.png"
$bc = $ws.Cells.Item(10, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(10, 3).Value = "合肥·星域动漫游戏嘉年华"
$ws.Cells.Item(10, 4).Value = "新站区东方大道288号 少荃体育中心"
$ws.Cells.Item(10, 5).Value = "2024.10.01 10:00-10.01 17:00"
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 58
$ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91878"
$ws.Cells.Item(10, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png"
$bc = $ws.Cells.Item(11, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(11, 3).Value = "合肥·第十五届次元之门动漫游戏博览会"
$ws.Cells.Item(11, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(11, 5).Value = "2024.10.01 09:30-10.02 17:30"
$ws.Cells.Item(11, 6).Value = 4543
$ws.Cells.Item(11, 7).Value = 70
$ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91133"
$ws.Cells.Item(11, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"
$bc = $ws.Cells.Item(12, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(12, 3).Value = "合肥·首届AT次元时代动漫游戏嘉年华"
$ws.Cells.Item(12, 4).Value = "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"
$ws.Cells.Item(12, 5).Value = "2024.10.01 09:30-10.03 17:00"
$ws.Cells.Item(12, 6).Value = 4391
$ws.Cells.Item(12, 7).Value = 68
$ws.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90908"
$ws.Cells.Item(12, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg"
$bc = $ws.Cells.Item(13, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-04"
$bc.Style = "Normal"
$ws.Cells.Item(13, 3).Value = "合肥·乐帮•崩原铁绝only同人首展"
$ws.Cells.Item(13, 4).Value = "丹霞路488号金星商业城三楼 迷鹿轰趴"
$ws.Cells.Item(13, 5).Value = "2024.10.04 10:00-10.05 16:30"
$ws.Cells.Item(13, 6).Value = 10
$ws.Cells.Item(13, 7).Value = 58
$ws.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91524"
$ws.Cells.Item(13, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png"
$ws.Range("A6").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Cells.Item(14, 1).Value = 13
$bc = $ws.Cells.Item(14, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-06"
$bc.Style = "Normal"
$ws.Cells.Item(14, 3).Value = "合肥·首届火影忍者同人only"
$ws.Cells.Item(14, 4).Value = "长江东路金太阳家具广场南门二楼 优极篮球馆"
$ws.Cells.Item(14, 5).Value = "2024.10.06 09:30-10.06 17:30"
$ws.Cells.Item(14, 6).Value = 18
$ws.Cells.Item(14, 7).Value = 75
$ws.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91658"
$ws.Cells.Item(14, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg"
$ws.Range("A6").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Cells.Item(15, 1).Value = 14
$bc = $ws.Cells.Item(15, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-26"
$bc.Style = "Normal"
$ws.Cells.Item(15, 3).Value = "合肥·W·A第五人格同人only2.0"
$ws.Cells.Item(15, 4).Value = "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)"
$ws.Cells.Item(15, 5).Value = "2024.10.26 09:30-10.26 17:00"
$ws.Cells.Item(15, 6).Value = 144
$ws.Cells.Item(15, 7).Value = 68
$ws.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91123"
$ws.Cells.Item(15, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 710
$ws.Cells.Item(3, 6).Value = 58
$ws.Cells.Item(4, 6).Value = 544
$ws.Cells.Item(5, 6).Value = 42
$bc = $ws.Cells.Item(7, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-15"
$bc.Style = "Normal"
$ws.Cells.Item(7, 3).Value = "蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）"
$ws.Cells.Item(7, 4).Value = "高新区望江西路888号 银泰百货（高新店）"
$ws.Cells.Item(7, 5).Value = "2024.09.15 10:00-10.02 22:00"
$ws.Cells.Item(7, 6).Value = 6
$ws.Cells.Item(7, 7).Value = 30
$ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91869"
$ws.Cells.Item(7, 9).Value = "//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png"
$bc = $ws.Cells.Item(8, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-16"
$bc.Style = "Normal"
$ws.Cells.Item(8, 3).Value = "肥西·星域动漫游戏嘉年华"
$ws.Cells.Item(8, 4).Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws.Cells.Item(8, 5).Value = "2024.09.16 10:00-09.16 17:00"
$ws.Cells.Item(8, 6).Value = 54
$ws.Cells.Item(8, 7).Value = 45
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws.Cells.Item(8, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"
$bc = $ws.Cells.Item(9, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-09-21"
$bc.Style = "Normal"
$ws.Cells.Item(9, 3).Value = "合肥·漫有引力动漫游戏嘉年华"
$ws.Cells.Item(9, 4).Value = "幸福路1号(筑梦集团·结婚产业园·B1幢) 费加罗宴会艺术中心(旗舰店)"
$ws.Cells.Item(9, 5).Value = "2024.09.21 10:00-09.21 17:00"
$ws.Cells.Item(9, 6).Value = 48
$ws.Cells.Item(9, 7).Value = 50
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90905"
$ws.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/gfeOndjM1723659151069.png"
$bc = $ws.Cells.Item(10, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(10, 3).Value = "合肥·星域动漫游戏嘉年华"
$ws.Cells.Item(10, 4).Value = "新站区东方大道288号 少荃体育中心"
$ws.Cells.Item(10, 5).Value = "2024.10.01 10:00-10.01 17:00"
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 58
$ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91878"
$ws.Cells.Item(10, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png"
$bc = $ws.Cells.Item(11, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(11, 3).Value = "合肥·第十五届次元之门动漫游戏博览会"
$ws.Cells.Item(11, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(11, 5).Value = "2024.10.01 09:30-10.02 17:30"
$ws.Cells.Item(11, 6).Value = 4543
$ws.Cells.Item(11, 7).Value = 70
$ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91133"
$ws.Cells.Item(11, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"
$bc = $ws.Cells.Item(12, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-01"
$bc.Style = "Normal"
$ws.Cells.Item(12, 3).Value = "合肥·首届AT次元时代动漫游戏嘉年华"
$ws.Cells.Item(12, 4).Value = "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"
$ws.Cells.Item(12, 5).Value = "2024.10.01 09:30-10.03 17:00"
$ws.Cells.Item(12, 6).Value = 4391
$ws.Cells.Item(12, 7).Value = 68
$ws.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90908"
$ws.Cells.Item(12, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg"
$bc = $ws.Cells.Item(13, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-04"
$bc.Style = "Normal"
$ws.Cells.Item(13, 3).Value = "合肥·乐帮•崩原铁绝only同人首展"
$ws.Cells.Item(13, 4).Value = "丹霞路488号金星商业城三楼 迷鹿轰趴"
$ws.Cells.Item(13, 5).Value = "2024.10.04 10:00-10.05 16:30"
$ws.Cells.Item(13, 6).Value = 10
$ws.Cells.Item(13, 7).Value = 58
$ws.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91524"
$ws.Cells.Item(13, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png"
$bc = $ws.Cells.Item(14, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-06"
$bc.Style = "Normal"
$ws.Cells.Item(14, 3).Value = "合肥·首届火影忍者同人only"
$ws.Cells.Item(14, 4).Value = "长江东路金太阳家具广场南门二楼 优极篮球馆"
$ws.Cells.Item(14, 5).Value = "2024.10.06 09:30-10.06 17:30"
$ws.Cells.Item(14, 6).Value = 18
$ws.Cells.Item(14, 7).Value = 75
$ws.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91658"
$ws.Cells.Item(14, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg"
$bc = $ws.Cells.Item(15, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-26"
$bc.Style = "Normal"
$ws.Cells.Item(15, 3).Value = "合肥·W·A第五人格同人only2.0"
$ws.Cells.Item(15, 4).Value = "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)"
$ws.Cells.Item(15, 5).Value = "2024.10.26 09:30-10.26 17:00"
$ws.Cells.Item(15, 6).Value = 144
$ws.Cells.Item(15, 7).Value = 68
$ws.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91123"
$ws.Cells.Item(15, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"
$bc = $ws.Cells.Item(16, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-10-26"
$bc.Style = "Normal"
$ws.Cells.Item(16, 3).Value = "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$ws.Cells.Item(16, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws.Cells.Item(16, 5).Value = "2024.10.26 19:30-10.26 21:00"
$ws.Cells.Item(16, 6).Value = 66
$ws.Cells.Item(16, 7).Value = 80
$ws.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90322"
$ws.Cells.Item(16, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"
$ws.Range("A6").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Cells.Item(17, 1).Value = 16
$bc = $ws.Cells.Item(17, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-11-09"
$bc.Style = "Normal"
$ws.Cells.Item(17, 3).Value = "合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会"
$ws.Cells.Item(17, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws.Cells.Item(17, 5).Value = "2024.11.09 19:30-11.09 21:00"
$ws.Cells.Item(17, 6).Value = 5
$ws.Cells.Item(17, 7).Value = 64
$ws.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90593"
$ws.Cells.Item(17, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg"
$ws.Range("A6").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Cells.Item(18, 1).Value = 17
$bc = $ws.Cells.Item(18, 2)
$bc.NumberFormat = "@"
$bc.Value = "2024-12-07"
$bc.Style = "Normal"
$ws.Cells.Item(18, 3).Value = "合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会"
$ws.Cells.Item(18, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws.Cells.Item(18, 5).Value = "2024.12.07 19:30-12.07 21:00"
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 56
$ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91608"
$ws.Cells.Item(18, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg"
